# The data rows (below the header) need to be re-sorted by the "Age"
# column (column C) in ascending order. Cell formatting (the green fill
# applied to the five oldest people) travels with the row data, exactly
# like Excel's native Range.Sort behavior.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the whole populated table, header included.
$dataRange = $ws.UsedRange

# Sort key: the "Age" column, first row of the range (used by Excel to
# figure out which column of the range to sort on).
$sortKey = $ws.Range("C1")

# Key1, Order1, Key2, Type2, Order2, Key3, Order3, Orientation,
# SortMethod, DataOption1, DataOption2, DataOption3 -- using the classic
# Range.Sort COM signature. Order1 = 1 (xlAscending), Orientation = 1
# (xlSortRows/xlTopToBottom), Header = xlYes (the range includes the
# header row, which must stay first and untouched).
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1, $true)
